# The dataprovider now supplies filepath + sheetname dynamically, so the
# second, now-unused sample sheet goes away and Sheet1 gains two more
# sample rows (one reusing the existing English/Amsterdam locale pair, one
# using a new Arabic/Canary locale pair) to exercise the new rows coming
# from the dataprovider.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$sheet2 = $wb.Worksheets.Item("Sheet2")
$null = $sheet2.Delete()

$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 now uses the Arabic (U.A.E.) / Atlantic-Canary locale pair instead
# of English (United States) / Europe-Amsterdam.
$ws.Range("E2").Value = "Arabic (U.A.E.)"
$ws.Range("F2").Value = "(UTC+00:00:00) Atlantic/Canary"

# Row 3: same shape as the original row 1/2, English/Amsterdam locale pair.
$ws.Range("A3").Value = "venkat"
$ws.Range("B3").Value = "p"
$ws.Range("C3").Value = "System - Administrator"
$ws.Range("D3").Value = "Administrator"
$ws.Range("E3").Value = "English (United States)"
$ws.Range("F3").Value = "(UTC+01:00:00) Europe/Amsterdam"

# Row 4: Arabic (U.A.E.) / Atlantic-Canary locale pair.
$ws.Range("A4").Value = "venkat"
$ws.Range("B4").Value = "p"
$ws.Range("C4").Value = "System - Administrator"
$ws.Range("D4").Value = "Administrator"
$ws.Range("E4").Value = "Arabic (U.A.E.)"
$ws.Range("F4").Value = "(UTC+00:00:00) Atlantic/Canary"

$null = $ws.Range("B8").Select()

Write-Host "done"
